# Fruta / hortaliza, semanal
# Insert one new weekly record at row 185 (pushing the existing rows 185-196
# down to 186-197) in the Chirimoya - Vega Modelo de Temuco sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 185; this shifts rows 185..196 down to 186..197
# and grows the used range from A1:T196 to A1:T197, matching existing
# row formatting (e.g. the date style on column D) via Excel's normal
# "insert inherits formatting from the row above" behaviour.
$ws.Rows.Item(185).Insert()

# Populate the new row with the new weekly observation.
$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 45166
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = "Fruta"
$ws.Range("G185").Value = 100107
$ws.Range("H185").Value = "Otros"
$ws.Range("I185").Value = 100107002
$ws.Range("J185").Value = "Chirimoya"
$ws.Range("K185").Value = "Cultivar IV Región"
$ws.Range("L185").Value = "Primera"
$ws.Range("M185").Value = 55
$ws.Range("N185").Value = 3000
$ws.Range("O185").Value = 3000
$ws.Range("P185").Value = 3000
$ws.Range("Q185").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R185").Value = "Provincia del Elquí"
$ws.Range("S185").Value = 3000
$ws.Range("T185").Value = 1
